# "proc enum in struct"
#
# 1) Primitives sheet: add a new "Primitive_ID" column (H) header.
#    (The used-range / dimension naturally grows to A1:H17 because rows
#    2-17 already have data in columns A-G; only the header cell is set.)
# 2) Types sheet: column N ("Alias") previously held a boolean-ish marker
#    (1 = "this field's IE_Type is a primitive", blank = "it is not").
#    It is now replaced by the actual numeric Primitive_ID that the field's
#    IE_Type resolves to in the Primitives sheet (or -1 when the IE_Type is
#    not a primitive / has no match there).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Primitives sheet
# ---------------------------------------------------------------------
$wsPrim = $wb.Worksheets.Item("Primitives")

$wsPrim.Cells.Item(1, 8).Value = "Primitive_ID"

# Match column H's width (14) to the sibling columns. Excel's ColumnWidth
# property is offset from the stored XML width by a small constant padding
# factor on this sheet's font; 13.17 round-trips to width="14" in the XML,
# same as the existing column G.
$wsPrim.Columns.Item(8).ColumnWidth = 13.17

# ---------------------------------------------------------------------
# Types sheet - column N ("Alias")
# ---------------------------------------------------------------------
$wsTypes = $wb.Worksheets.Item("Types")

$aliasValues = @(
    @{Row=2;  Value=13},
    @{Row=3;  Value=13},
    @{Row=4;  Value=13},
    @{Row=5;  Value=13},
    @{Row=6;  Value=13},
    @{Row=7;  Value=13},
    @{Row=8;  Value=-1},
    @{Row=9;  Value=9},
    @{Row=10; Value=13},
    @{Row=11; Value=-1},
    @{Row=12; Value=-1},
    @{Row=13; Value=-1},
    @{Row=14; Value=-1},
    @{Row=15; Value=-1},
    @{Row=16; Value=-1},
    @{Row=17; Value=-1},
    @{Row=18; Value=-1},
    @{Row=19; Value=-1},
    @{Row=20; Value=-1},
    @{Row=21; Value=-1},
    @{Row=22; Value=-1},
    @{Row=23; Value=-1},
    @{Row=24; Value=6},
    @{Row=25; Value=6},
    @{Row=26; Value=10},
    @{Row=27; Value=-1},
    @{Row=28; Value=6},
    @{Row=29; Value=-1},
    @{Row=30; Value=-1},
    @{Row=31; Value=-1},
    @{Row=32; Value=-1},
    @{Row=33; Value=5},
    @{Row=34; Value=-1},
    @{Row=35; Value=-1},
    @{Row=36; Value=-1},
    @{Row=37; Value=-1},
    @{Row=38; Value=-1},
    @{Row=39; Value=-1},
    @{Row=40; Value=-1},
    @{Row=41; Value=-1},
    @{Row=42; Value=-1},
    @{Row=43; Value=-1},
    @{Row=44; Value=-1},
    @{Row=45; Value=-1},
    @{Row=46; Value=-1},
    @{Row=47; Value=8},
    @{Row=48; Value=-1},
    @{Row=49; Value=-1},
    @{Row=50; Value=-1},
    @{Row=51; Value=8},
    @{Row=52; Value=-1},
    @{Row=53; Value=8},
    @{Row=54; Value=-1},
    @{Row=55; Value=8},
    @{Row=56; Value=-1},
    @{Row=57; Value=8},
    @{Row=58; Value=-1},
    @{Row=59; Value=6},
    @{Row=60; Value=6},
    @{Row=61; Value=-1},
    @{Row=62; Value=6},
    @{Row=63; Value=-1},
    @{Row=64; Value=-1},
    @{Row=65; Value=-1},
    @{Row=66; Value=-1}
)

foreach ($entry in $aliasValues) {
    $wsTypes.Cells.Item($entry.Row, 14).Value = $entry.Value
}
